# Weekly update: a new daily price record is inserted at row 62, pushing the
# existing records (rows 62-128) down by one row (to rows 63-129).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 62; Excel shifts rows 62:128 down to 63:129
# and the sheet dimension grows from A1:R128 to A1:R129 automatically.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A62").Value = 4
$ws.Range("B62").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C62").Value = "Los Lagos"
$ws.Range("D62").Value = 44778
$ws.Range("E62").Value = 10
$ws.Range("F62").Value = 100112022
$ws.Range("G62").Value = "Arveja Verde"
$ws.Range("H62").Value = "Perfection"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 45
$ws.Range("K62").Value = 45000
$ws.Range("L62").Value = 45000
$ws.Range("M62").Value = 45000
$ws.Range("N62").Value = "$/malla 25 kilos"
$ws.Range("O62").Value = "Provincia de Huasco"
$ws.Range("P62").Value = 1800
$ws.Range("Q62").Value = 25
$ws.Range("R62").Value = "Hortaliza"
